$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1298.7831
$ws.Range("I15").Value = 1298.7831
$ws.Range("K15").Value = 3896.3493
$ws.Range("M15").Value = -3727.3493

# Row 43
$ws.Range("H43").Value = 1176.6
$ws.Range("I43").Value = 1015.26086
$ws.Range("J43").Value = 1485.8334
$ws.Range("K43").Value = 1015.26086
$ws.Range("L43").Value = 1485.8334
$ws.Range("M43").Value = -946.26086
$ws.Range("N43").Value = -1623.8334

# Row 138
$ws.Range("H138").Value = 4164.772
$ws.Range("I138").Value = 1585.0588
$ws.Range("J138").Value = 7978.2607
$ws.Range("K138").Value = 4755.1764
$ws.Range("L138").Value = 23934.7821
$ws.Range("M138").Value = 384.8235999999997
$ws.Range("N138").Value = -34214.7821

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 13374133
$ws.Range("I45").Value = 23881650
$ws.Range("J45").Value = 930.4545000000001
$ws.Range("K45").Value = 23881650
$ws.Range("L45").Value = 930.4545000000001
$ws.Range("M45").Value = -23881273
$ws.Range("N45").Value = -1684.4545

# Row 61
$ws.Range("H61").Value = 2123.8572
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 2123.8572
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 2123.8572
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -2547.8572

# Row 74
$ws.Range("H74").Value = 1650.3658
$ws.Range("I74").Value = 1307.4783
$ws.Range("J74").Value = 2088.5
$ws.Range("K74").Value = 1307.4783
$ws.Range("L74").Value = 2088.5
$ws.Range("M74").Value = -433.4783
$ws.Range("N74").Value = -3836.5

# Row 77
$ws.Range("H77").Value = 1650.3658
$ws.Range("I77").Value = 1307.4783
$ws.Range("J77").Value = 2088.5
$ws.Range("K77").Value = 6537.3915
$ws.Range("L77").Value = 10442.5
$ws.Range("M77").Value = -2169.3915
$ws.Range("N77").Value = -19178.5

# Row 97
$ws.Range("H97").Value = 1394.4117
$ws.Range("I97").Value = 1194
$ws.Range("J97").Value = 2897.5
$ws.Range("K97").Value = 1194
$ws.Range("L97").Value = 2897.5
$ws.Range("M97").Value = -698
$ws.Range("N97").Value = -3889.5

# Row 102
$ws.Range("H102").Value = 1985.5834
$ws.Range("I102").Value = 1859.9524
$ws.Range("J102").Value = 2865
$ws.Range("K102").Value = 1859.9524
$ws.Range("L102").Value = 2865
$ws.Range("M102").Value = -237.9523999999999
$ws.Range("N102").Value = -6109

# Row 122
$ws.Range("H122").Value = 1476.8
$ws.Range("I122").Value = 1486.7059
$ws.Range("J122").Value = 1140
$ws.Range("K122").Value = 4460.1177
$ws.Range("L122").Value = 3420
$ws.Range("M122").Value = -2010.1177
$ws.Range("N122").Value = -8320

# Row 136
$ws.Range("H136").Value = 2123.8572
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 2123.8572
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 6371.571599999999
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -11471.5716

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2005.8422
$ws.Range("I86").Value = 2208.077
$ws.Range("J86").Value = 1567.6666
$ws.Range("K86").Value = 2208.077
$ws.Range("L86").Value = 1567.6666
$ws.Range("M86").Value = -1085.077
$ws.Range("N86").Value = -3813.6666

# Row 89
$ws.Range("H89").Value = 2005.8422
$ws.Range("I89").Value = 2208.077
$ws.Range("J89").Value = 1567.6666
$ws.Range("K89").Value = 11040.385
$ws.Range("L89").Value = 7838.333000000001
$ws.Range("M89").Value = -5424.385000000002
$ws.Range("N89").Value = -19070.333

# Row 94
$ws.Range("H94").Value = 557.3
$ws.Range("I94").Value = 508.1111
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 508.1111
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -57.11110000000002
$ws.Range("N94").Value = -1902

# Row 132
$ws.Range("H132").Value = 42491.11
$ws.Range("J132").Value = 42491.11
$ws.Range("L132").Value = 42491.11
$ws.Range("N132").Value = -52611.11

$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 459.5
$ws.Range("I107").Value = 435.6316
$ws.Range("J107").Value = 913
$ws.Range("K107").Value = 435.6316
$ws.Range("L107").Value = 913
$ws.Range("M107").Value = 1484.3684
$ws.Range("N107").Value = -4753

# Row 134
$ws.Range("H134").Value = 12821584
$ws.Range("I134").Value = 831.85297
$ws.Range("K134").Value = 2495.55891
$ws.Range("M134").Value = 39.4410899999998

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 7146339.5
$ws.Range("I122").Value = 25000170
$ws.Range("K122").Value = 225001530
$ws.Range("M122").Value = -224999080

# Row 131
$ws.Range("H131").Value = 854.04126
$ws.Range("J131").Value = 878.5730600000001
$ws.Range("L131").Value = 2635.71918
$ws.Range("N131").Value = -12715.71918

# Row 141
$ws.Range("H141").Value = 3015.75
$ws.Range("I141").Value = 3010
$ws.Range("K141").Value = 9030
$ws.Range("M141").Value = -3850

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2386.2144
$ws.Range("I102").Value = 1932.4546
$ws.Range("K102").Value = 1932.4546
$ws.Range("M102").Value = -310.4546

# Row 122
$ws.Range("H122").Value = 19294514
$ws.Range("I122").Value = 27011094
$ws.Range("J122").Value = 3062.5
$ws.Range("K122").Value = 81033282
$ws.Range("L122").Value = 9187.5
$ws.Range("M122").Value = -81030832
$ws.Range("N122").Value = -14087.5

# Row 126
$ws.Range("H126").Value = 4507.467
$ws.Range("J126").Value = 4645.4546
$ws.Range("L126").Value = 13936.3638
$ws.Range("N126").Value = -18876.3638

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1597.5555
$ws.Range("I7").Value = 1297.25
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 1297.25
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -1185.25
$ws.Range("N7").Value = -4224

# Row 16
$ws.Range("H16").Value = 721.5454999999999
$ws.Range("I16").Value = 721.5454999999999
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 721.5454999999999
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -551.5454999999999
$ws.Range("N16").ClearContents()

# Row 74
$ws.Range("H74").Value = 18000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 18000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 18000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -19996

# Row 77
$ws.Range("H77").Value = 18000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 18000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 54000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -63984

# Row 93
$ws.Range("H93").Value = 1132
$ws.Range("I93").Value = 1230.2858
$ws.Range("J93").Value = 1055.5555
$ws.Range("K93").Value = 1230.2858
$ws.Range("L93").Value = 1055.5555
$ws.Range("M93").Value = 17.71419999999989
$ws.Range("N93").Value = -3551.5555

# Row 100
$ws.Range("H100").Value = 6899.5
$ws.Range("I100").Value = 10000
$ws.Range("J100").Value = 3799
$ws.Range("K100").Value = 10000
$ws.Range("L100").Value = 3799
$ws.Range("M100").Value = -9459
$ws.Range("N100").Value = -4881

# Row 106
$ws.Range("H106").Value = 333366660
$ws.Range("J106").Value = 333366660
$ws.Range("L106").Value = 333366660
$ws.Range("N106").Value = -333369184

# Row 126
$ws.Range("H126").Value = 1597.5555
$ws.Range("I126").Value = 1297.25
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 3891.75
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -1421.75
$ws.Range("N126").Value = -16940

# Row 136
$ws.Range("H136").Value = 93954400
$ws.Range("I136").Value = 123812510
$ws.Range("J136").Value = 58827210
$ws.Range("K136").Value = 371437530
$ws.Range("L136").Value = 176481630
$ws.Range("M136").Value = -371434980
$ws.Range("N136").Value = -176486730

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 3789.1428
$ws.Range("I62").Value = 3694.8
$ws.Range("J62").Value = 4025
$ws.Range("K62").Value = 3694.8
$ws.Range("L62").Value = 4025
$ws.Range("M62").Value = -3070.8
$ws.Range("N62").Value = -5273

# Row 65
$ws.Range("H65").Value = 3789.1428
$ws.Range("I65").Value = 3694.8
$ws.Range("J65").Value = 4025
$ws.Range("K65").Value = 18474
$ws.Range("L65").Value = 20125
$ws.Range("M65").Value = -15354
$ws.Range("N65").Value = -26365

# Row 81
$ws.Range("H81").Value = 817.5
$ws.Range("I81").Value = 542
$ws.Range("J81").Value = 1093
$ws.Range("K81").Value = 1084
$ws.Range("L81").Value = 2186
$ws.Range("M81").Value = -23
$ws.Range("N81").Value = -4308

# Row 84
$ws.Range("H84").Value = 817.5
$ws.Range("I84").Value = 542
$ws.Range("J84").Value = 1093
$ws.Range("K84").Value = 5420
$ws.Range("L84").Value = 10930
$ws.Range("M84").Value = -116
$ws.Range("N84").Value = -21538

# Row 100
$ws.Range("H100").Value = 509.57895
$ws.Range("I100").Value = 431.73334
$ws.Range("J100").Value = 801.5
$ws.Range("K100").Value = 863.46668
$ws.Range("L100").Value = 1603
$ws.Range("M100").Value = -322.46668
$ws.Range("N100").Value = -2685

# Row 107
$ws.Range("H107").Value = 313.4
$ws.Range("I107").Value = 365.5
$ws.Range("J107").Value = 105
$ws.Range("K107").Value = 1096.5
$ws.Range("L107").Value = 315
$ws.Range("M107").Value = 823.5
$ws.Range("N107").Value = -4155

# Row 113
$ws.Range("H113").Value = 71432480
$ws.Range("I113").Value = 90913656
$ws.Range("K113").Value = 272740968
$ws.Range("M113").Value = -272738798

# Row 122
$ws.Range("H122").Value = 65250.5
$ws.Range("I122").Value = 101700.8
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 305102.4
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -302652.4
$ws.Range("N122").Value = -18400

# Row 136
$ws.Range("H136").Value = 2586.913
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 2586.913
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 7760.739
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -12860.739
